$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.080.88"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.660.55"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "'208.18"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").Value = "'0.5177"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "'0.2582"
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("D9").Value = "'0.06293"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "'20.92"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").Value = "'0.07538"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "1.665.66"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").Value = "'4.403"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "'0.5390"
$ws.Range("E14").Value = "  -3.98%  "
$ws.Range("D15").Value = "'66.12"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "0.0₅7928"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "26.089.74"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "'4.696"
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("D20").Value = "'187.81"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("E21").Value = "  -2.63%  "
$ws.Range("D22").Value = "'6.189"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "'147.90"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").Value = "'0.1207"
$ws.Range("E25").Value = "  -3.62%  "
$ws.Range("D26").Value = "'7.384"
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("D27").Value = "'15.66"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "'1.381"
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("D29").Value = "'0.06039"
$ws.Range("E29").Value = "  -5.32%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").Value = "'3.470"
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").Value = "'3.397"
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("D33").Value = "'1.634"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").Value = "'0.9840"
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D35").Value = "'2.387"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").Value = "'2.752"
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("D37").Value = "'0.5881"
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("D38").Value = "1.104.82"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").Value = "'0.01593"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("D40").Value = "'5.953"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("D41").Value = "'0.8486"
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "'99.88"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "1.812.78"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").Value = "0.0₈109"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").Value = "'55.10"
$ws.Range("E46").Value = "  -2.65%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "'8.026"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "'0.05227"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").Value = "'0.4240"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "'5.857"
$ws.Range("E51").Value = "  -1.15%  "
